$d = $word.ActiveDocument
$p11 = $d.Paragraphs.Item(11)
$p11.Range.Text = "Group Website"
$p11 = $d.Paragraphs.Item(11)
$r = $d.Range($p11.Range.Start, $p11.Range.End - 1)
$r.Collapse(0)
$r.InsertAfter("`r`rhttps://s3728065.github.io/A2-24/`r`rPrepared by:")

# Now try to find the empty paragraph (#12) and see its range precisely
$p12 = $d.Paragraphs.Item(12)
Write-Output "p12: $($p12.Range.Start)-$($p12.Range.End) text=[$($p12.Range.Text)]"
$p12.Range.Font.Bold
